# Scrum sprint 2 fertig gestellt
#
# Updates the "Sprint Backlog" sheet:
#   - corrects a few Sprint-1/2 rows (priority High -> Medium, status -> Moved/Done)
#   - adds the missing "Design verbessern" story (row 17, id 2.9)
#   - fills in the whole Sprint 3 backlog (rows 18-24, were empty stub rows before)
#
# NOTE on ordering: new shared-string entries are appended to the workbook's
# string table in first-use order. The cell writes below are deliberately
# sequenced (first-touch of each brand-new label) to reproduce that order:
#   Design verbessern, Medium, Navigation von Appointment zu Patient View,
#   Abschlussdialog, "Testing ", Unit Tests, Testdaten generiern,
#   Design responsive machen, UI, Low

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")

$arrow = [char]8594   # "→"

# 1) "Design verbessern" must become the first brand-new shared string.
$ws.Cells.Item(17, 3).Value  = "Design verbessern"

# 2) "Medium" becomes the second brand-new shared string (priority fixes).
$ws.Cells.Item(5, 8).Value   = "Medium"   # H5
$ws.Cells.Item(6, 8).Value   = "Medium"   # H6
$ws.Cells.Item(7, 8).Value   = "Medium"   # H7
$ws.Cells.Item(9, 8).Value   = "Medium"   # H9
$ws.Cells.Item(10, 8).Value  = "Medium"   # H10
$ws.Cells.Item(14, 8).Value  = "Medium"   # H14
$ws.Cells.Item(17, 8).Value  = "Medium"   # H17
$ws.Cells.Item(19, 8).Value  = "Medium"   # H19

# 3) "Navigation von Appointment zu Patient View"
$ws.Cells.Item(20, 3).Value  = "Navigation von Appointment zu Patient View"

# 4) "Abschlussdialog"
$ws.Cells.Item(21, 3).Value  = "Abschlussdialog"

# 5) "Testing " (note trailing space) must be created before "Unit Tests"
$ws.Cells.Item(22, 5).Value  = "Testing "
$ws.Cells.Item(23, 5).Value  = "Testing "

# 6) "Unit Tests"
$ws.Cells.Item(22, 3).Value  = "Unit Tests"

# 7) "Testdaten generiern"
$ws.Cells.Item(23, 3).Value  = "Testdaten generiern"

# 8) "Design responsive machen"
$ws.Cells.Item(24, 3).Value  = "Design responsive machen"

# 9) "UI" (plain, distinct from the existing "UI, Controller" label)
$ws.Cells.Item(17, 5).Value  = "UI"
$ws.Cells.Item(24, 5).Value  = "UI"

# 10) "Low"
$ws.Cells.Item(24, 8).Value  = "Low"

# --- Status corrections on existing Sprint 1/2 rows -----------------------
$ws.Cells.Item(11, 12).Value = "Moved"   # L11
$ws.Cells.Item(14, 11).Value = 7         # K14 (remaining effort)
$ws.Cells.Item(14, 12).Value = "Moved"   # L14
$ws.Cells.Item(15, 12).Value = "Done"    # L15

# --- Remaining (non string-table-affecting) cells of row 17 "Design verbessern" ---
$ws.Cells.Item(17, 1).Value  = 2.9
$ws.Cells.Item(17, 2).Value  = 2
$ws.Cells.Item(17, 6).Value  = "Hannes"
$ws.Cells.Item(17, 7).Value  = "Nic"
$ws.Cells.Item(17, 9).Value  = 7
$ws.Cells.Item(17, 10).Value = 7
$ws.Cells.Item(17, 11).Value = 7
$ws.Cells.Item(17, 12).Value = "Done"

# --- Sprint 3 backlog, rows 18-24 (remaining columns) ----------------------

# 3.1 Patient Information view -> Tab Todo
$ws.Cells.Item(18, 1).Value  = 3.1
$ws.Cells.Item(18, 2).Value  = 3
$ws.Cells.Item(18, 3).Value  = "Patient Information view " + $arrow + " Tab Todo"
$ws.Cells.Item(18, 5).Value  = "UI, Controller"
$ws.Cells.Item(18, 6).Value  = "Ken"
$ws.Cells.Item(18, 7).Value  = "Joel"
$ws.Cells.Item(18, 8).Value  = "High"
$ws.Cells.Item(18, 9).Value  = 10
$ws.Cells.Item(18, 10).Value = 14
$ws.Cells.Item(18, 11).Value = 10
$ws.Cells.Item(18, 12).Value = "In Progress"

# 3.2 Patient Information view -> Tab PatientInformation
$ws.Cells.Item(19, 1).Value  = 3.2
$ws.Cells.Item(19, 2).Value  = 3
$ws.Cells.Item(19, 3).Value  = "Patient Information view " + $arrow + " Tab PatientInformation"
$ws.Cells.Item(19, 5).Value  = "UI, Controller"
$ws.Cells.Item(19, 6).Value  = "Hannes"
$ws.Cells.Item(19, 7).Value  = "Nic"
$ws.Cells.Item(19, 9).Value  = 14
$ws.Cells.Item(19, 10).Value = 14
$ws.Cells.Item(19, 11).Value = 10
$ws.Cells.Item(19, 12).Value = "In Progress"

# 3.3 Navigation von Appointment zu Patient View
$ws.Cells.Item(20, 1).Value  = 3.3
$ws.Cells.Item(20, 2).Value  = 3
$ws.Cells.Item(20, 5).Value  = "UI, Controller"
$ws.Cells.Item(20, 6).Value  = "Nic"
$ws.Cells.Item(20, 7).Value  = "Hannes"
$ws.Cells.Item(20, 8).Value  = "High"
$ws.Cells.Item(20, 9).Value  = 7
$ws.Cells.Item(20, 10).Value = 7
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = "In Progress"

# 3.4 Abschlussdialog
$ws.Cells.Item(21, 1).Value  = 3.4
$ws.Cells.Item(21, 2).Value  = 3
$ws.Cells.Item(21, 5).Value  = "UI, Controller"
$ws.Cells.Item(21, 6).Value  = "Gabor"
$ws.Cells.Item(21, 7).Value  = "Quentin"
$ws.Cells.Item(21, 8).Value  = "High"
$ws.Cells.Item(21, 9).Value  = 7
$ws.Cells.Item(21, 10).Value = 7
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = "In Progress"

# 3.5 Unit Tests
$ws.Cells.Item(22, 1).Value  = 3.5
$ws.Cells.Item(22, 2).Value  = 3
$ws.Cells.Item(22, 6).Value  = "Quentin"
$ws.Cells.Item(22, 7).Value  = "Gabor"
$ws.Cells.Item(22, 8).Value  = "High"
$ws.Cells.Item(22, 9).Value  = 7
$ws.Cells.Item(22, 10).Value = 7
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = "In Progress"

# 3.6 Testdaten generiern
$ws.Cells.Item(23, 1).Value  = 3.6
$ws.Cells.Item(23, 2).Value  = 3
$ws.Cells.Item(23, 6).Value  = "Joel"
$ws.Cells.Item(23, 7).Value  = "Ken"
$ws.Cells.Item(23, 8).Value  = "High"
$ws.Cells.Item(23, 9).Value  = 7
$ws.Cells.Item(23, 10).Value = 7
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = "In Progress"

# 3.7 Design responsive machen
$ws.Cells.Item(24, 1).Value  = 3.7
$ws.Cells.Item(24, 2).Value  = 3
$ws.Cells.Item(24, 6).Value  = "Hannes"
$ws.Cells.Item(24, 7).Value  = "Nic"
$ws.Cells.Item(24, 9).Value  = 3
$ws.Cells.Item(24, 10).Value = 3
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = "In Progress"

# --- Column C width (widened so the long new story names fit) -------------
$ws.Columns.Item(3).ColumnWidth = 39.5

# --- View state: scroll down one row, select D16 (matches author's edit) --
$ws.Activate()
$ws.Range("D16").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
